$d = $word.ActiveDocument

# --- Locate the paragraph that currently ends the "how to access the server"
#     section, i.e. the one whose text is the ClaimCenter.do URL line. We add
#     the new "default username..." paragraph right after it. ---
$urlParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*ClaimCenter.do*") {
        $urlParaIndex = $i
        break
    }
}

# Insert a brand new paragraph right after it.
$d.Paragraphs.Item($urlParaIndex).Range.InsertParagraphAfter() | Out-Null
$newParaIndex = $urlParaIndex + 1

# Give the new paragraph its text. A trailing placeholder character is used
# so the "_GoBack" bookmark (which Word always keeps collapsed at the very
# end of a run) can be anchored precisely after the real text, then the
# placeholder is removed, leaving bookmarkStart/bookmarkEnd directly after
# the run as Word itself would produce.
$d.Paragraphs.Item($newParaIndex).Range.Text = "default username is su and password is gwX"

$newRange = $d.Paragraphs.Item($newParaIndex).Range
$placeholder = $d.Range($newRange.End - 2, $newRange.End - 1)

# Adding a bookmark named "_GoBack" automatically removes any previous
# bookmark of that name elsewhere in the document (bookmark names are
# unique), so this both relocates it and drops the stale one near the
# "https://www.microsoft.com/..." paragraph in one step.
$d.Bookmarks.Add("_GoBack", $placeholder) | Out-Null

# Drop the placeholder character now that the bookmark is anchored.
$newRange2 = $d.Paragraphs.Item($newParaIndex).Range
$d.Range($newRange2.End - 2, $newRange2.End - 1).Delete()
